# RPA datasets push 2024-03-20
# Insert a new IPO row ("제일엠앤에스(구.제일기공)") into the 7th data row of the
# "02_38커뮤니케이션(최근일자기준)" sheet, push the remaining rows down by one,
# and drop the oldest row that falls off the bottom of the tracked window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Shift existing row 7 (and everything below it) down by one row.
$ws.Rows.Item(7).Insert()

# Populate the newly-inserted row 7 with the new IPO's data.
$ws.Cells.Item(7, 1).Value = "제일엠앤에스(구.제일기공)"
$ws.Cells.Item(7, 2).Value = "2024.04.05~04.12"
$ws.Cells.Item(7, 3).Value = "15,000~18,000"
$ws.Cells.Item(7, 4).Value = "-"
$ws.Cells.Item(7, 5).Value = 36000
$ws.Cells.Item(7, 6).Value = "케이비증권"

# The table keeps a fixed 20-row window, so drop the row that was pushed
# past the end (originally row 21, now shifted to row 22).
$ws.Rows.Item(22).Delete()
